$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E: hide it (was visible) ---
$ws.Columns("E").Hidden = $true

# --- Row 17 (Date 2020-02-12): fill in Task# / Task / Sub-activity / Break hours ---
$ws.Range("B17").Value = 13
$ws.Range("C17").Value = "Break"
$ws.Range("D17").Value = "Break"
$ws.Range("F17").Value = 0

# --- Row 18 (Date 2020-02-13): fill in Task# / Task / Sub-activity / Break hours ---
$ws.Range("B18").Value = 14
$ws.Range("C18").Value = "Break"
$ws.Range("D18").Value = "Break"
$ws.Range("F18").Value = 0

# --- Row 19 (Date 2020-02-14): new task entry ---
$ws.Range("B19").Value = 15
$ws.Range("C19").Value = "Create nested XSD Tags as hierarchical Checkboxes"
$ws.Range("C19").Font.Bold = $true
$ws.Range("C19").HorizontalAlignment = -4108
$ws.Range("C19").VerticalAlignment = -4108
$ws.Range("C19").WrapText = $true
$ws.Range("C19").Borders.LineStyle = 1

$ws.Range("D19").Value = "Tried to modify the DFS approach to filter out non-element tags"
$ws.Range("D19").Characters(52, 8).Font.Italic = $true
$ws.Range("D19").Characters(60, 4).Font.Italic = $false

$ws.Range("E19").Value = 0.1
$ws.Range("E19").NumberFormat = "0%"
$ws.Range("F19").Value = 3
$ws.Rows(19).RowHeight = 29

# --- Row 20 (Date 2020-02-15, Saturday): new task entry + weekend highlight ---
$ws.Range("A20").Interior.Color = 65535
$ws.Range("B20").Value = 16
$ws.Range("C20").Value = "Nested XSD Tags & Final Dataframe writing to CSV"
$ws.Range("D20").Value = "Successfully displayed nested hierarchical set of XSD Tags.`nWorked on creating final dataframe logic and successfully managed to create Pyspark config.csv file as per discussion with Pooja."
$ws.Range("E20").Value = 0.7
$ws.Range("E20").NumberFormat = "0%"
$ws.Range("F20").Value = 6
$ws.Rows(20).RowHeight = 43.5

# --- Selection / viewport ---
$ws.Range("F13").Select()

Write-Output "edit applied"
